$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns("A").ColumnWidth = 27.6640625
